$wb = $excel.ActiveWorkbook

# --- Sheet1 (PERMANOVA) ---
$ws1 = $wb.Worksheets.Item("PERMANOVA")

# Row 2
$ws1.Range("B2").Value = 0.8740278822251073
$ws1.Range("C2").Value = 0.06870153376771632
$ws1.Range("D2").Value = 2.922387458680074
$ws1.Range("E2").Value = 0.01346

# Row 3
$ws1.Range("A3").Value = 9
$ws1.Range("B3").Value = 5.268311750865272
$ws1.Range("C3").Value = 0.4141070382440141
$ws1.Range("D3").Value = 1.957228202676015
$ws1.Range("E3").Value = 0.00132

# Row 4 (new row, shifted from old row3's A4=32 data being pushed down)
$ws1.Range("A4").Value = 22
$ws1.Range("B4").Value = 6.579761814895401
$ws1.Range("C4").Value = 0.51719142798827

# Row 5 (new row, previously row4 data: A=32,B=12.46...,C=1)
$ws1.Range("A5").Value = 32
$ws1.Range("B5").Value = 12.72210144798578
$ws1.Range("C5").Value = 1

# --- Sheet2 (PERMDISP) ---
$ws2 = $wb.Worksheets.Item("PERMDISP")

$ws2.Range("B2").Value = 0.00002994149613517597
$ws2.Range("C2").Value = 0.00002994149613517597
$ws2.Range("D2").Value = 0.001278238776438813
$ws2.Range("F2").Value = 0.98

$ws2.Range("B3").Value = 0.7261447526857168
$ws2.Range("C3").Value = 0.02342402428018441
